$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string text in B2 (row for the new "PROPAG" resource description),
# replacing the old "RECURSOS DECORRENTES DA DESVINCULAÇÃO..." entry.
$ws.Range("B2").Value = "RECURSOS DO ESTADO APLICADOS NOS INVESTIMENTOS RELATIVOS À LC 212/2025 – PROPAG"

# Move the active selection to B17 to match the saved view state.
$ws.Range("B17").Select()
